# Apply the edit described by the diff:
#  - Insert two new data rows right before the current row 144
#    (shifting the existing rows 144..239 down to 146..241).
#  - Populate the two new rows with their data.
#
# Both new rows share the constant columns used throughout the sheet
# (A,B,C,E,F,G,N,Q,R) and carry new values for D,H,I,J,K,L,M,O,P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 144, pushing current 144:241 -> 146:241 (well, 144:239 -> 146:241)
$ws.Rows("144:145").Insert()

# Common / constant values for every data row in this sheet
$colA = 7
$colB = "Terminal Hortofrutícola Agro Chillán"
$colC = "Ñuble"
$colE = 16
$colF = 100112045
$colG = "Zapallo"
$colN = '$/kilo (volumen en unidades)'
$colQ = 1
$colR = "Hortaliza"

# ---- New row 144 ----
$ws.Range("A144").Value = $colA
$ws.Range("B144").Value = $colB
$ws.Range("C144").Value = $colC
$ws.Range("D144").Value = 44981
$ws.Range("E144").Value = $colE
$ws.Range("F144").Value = $colF
$ws.Range("G144").Value = $colG
$ws.Range("H144").Value = "Camote"
$ws.Range("I144").Value = "1a (cosecha)"
$ws.Range("J144").Value = 400
$ws.Range("K144").Value = 400
$ws.Range("L144").Value = 450
$ws.Range("M144").Value = 425
$ws.Range("N144").Value = $colN
$ws.Range("O144").Value = "Región del Maule"
$ws.Range("P144").Value = 425
$ws.Range("Q144").Value = $colQ
$ws.Range("R144").Value = $colR

# ---- New row 145 ----
$ws.Range("A145").Value = $colA
$ws.Range("B145").Value = $colB
$ws.Range("C145").Value = $colC
$ws.Range("D145").Value = 44981
$ws.Range("E145").Value = $colE
$ws.Range("F145").Value = $colF
$ws.Range("G145").Value = $colG
$ws.Range("H145").Value = "Camote"
$ws.Range("I145").Value = "2a (cosecha)"
$ws.Range("J145").Value = 400
$ws.Range("K145").Value = 300
$ws.Range("L145").Value = 350
$ws.Range("M145").Value = 325
$ws.Range("N145").Value = $colN
$ws.Range("O145").Value = "Región del Maule"
$ws.Range("P145").Value = 325
$ws.Range("Q145").Value = $colQ
$ws.Range("R145").Value = $colR
